# Removed Test Case Inter-Dependency
#
# The "Input" sheet's row 5 (externalid / 4404) represented a hard-coded
# foreign-key value from a previous, unrelated test run. It is removed so
# this test case no longer depends on that prior run; all subsequent rows
# shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Delete the whole row (shifts rows 6:11 up to 5:10, and updates the
# sheet's used-range dimension accordingly).
$ws.Rows(5).EntireRow.Delete()

# Reselect the last populated cell, matching the post-edit cursor position.
$ws.Activate()
$ws.Range("B10").Select()
